$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("username|string" column), shifting
# the existing username/password/nickname/email/avatar columns one to the
# right.
$ws.Range("C1").EntireColumn.Insert()

# New column header + data: isAI|bool
$ws.Range("C1").Value = "isAI|bool"
$ws.Range("C3").Value = 1

# Size the new column similarly to the other (best-fit) columns.
$ws.Range("C1:C3").ColumnWidth = 8.45

# The hyperlink on the email cell doesn't automatically rebase when the
# column is inserted, so re-anchor it from F2 (old position) to G2 (new
# position). Re-apply the Hyperlink cell style afterwards since .Add()
# creates a fresh style rather than reusing the existing one.
$ws.Range("F2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:buptforeverbean@gmail.com")
$ws.Range("G2").Style = "Hyperlink"

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("F8").Select() | Out-Null
